{"js": "const body = context.document.body;\nconst results = body.search(\"Sole, on-going DevOps role. \", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target phrase not found\");\n}\n\nresults.items[0].insertText(\"Continued to own DevOps role. \", \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Sole, on-going DevOps role. \"\n$find.Replacement.Text = \"Continued to own DevOps role. \"\n$find.Forward = $true\n$find.Wrap = 1\n$find.MatchCase = $true\n$find.Execute($find.Text, $find.MatchCase, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
